$d = $word.ActiveDocument

# Remove the placeholder run "vnpt.SiteAddress" that follows "Địa chỉ: "
# leaving the label text itself untouched.
$d.Content.Find.Execute("vnpt.SiteAddress", $false, $false, $false, $false,
                         $false, $true, 1, $false, "", 2)
